$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.945.81'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '2.787.80'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '357.94'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.70'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.560'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.82%  '
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("E10").Value = '  -2.65%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.134'
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '0.0848'
$ws.Range("E12").Value = '  -1.84%  '
$ws.Range("D13").Value = '19.49'
$ws.Range("E13").Value = '  -3.76%  '
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("D15").Value = '3.227.33'
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '2.794.80'
$ws.Range("E16").Value = '  -2.12%  '
$ws.Range("D17").Value = '0.948'
$ws.Range("E17").Value = '  +1.97%  '
$ws.Range("D18").Value = '51.870.40'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").Value = '13.17'
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("D22").Value = '0.0₃0977'
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("D23").Value = '70.28'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '270.47'
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("E25").Value = '  -4.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.50'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  +16.08%  '
$ws.Range("E29").Value = '  -0.73%  '
$ws.Range("D30").Value = '2.14'
$ws.Range("E30").Value = '  -5.00%  '
$ws.Range("D31").Value = '0.0469'
$ws.Range("E31").Value = '  -3.43%  '
$ws.Range("D32").Value = '52.14'
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("D33").Value = '34.73'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").Value = '5.75'
$ws.Range("E34").Value = '  -2.74%  '
$ws.Range("D35").Value = '0.0845'
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").Value = '5.18'
$ws.Range("E36").Value = '  -6.17%  '
$ws.Range("D38").Value = '18.87'
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("E39").Value = '  -2.90%  '
$ws.Range("D40").Value = '1.99'
$ws.Range("E40").Value = '  -4.12%  '
$ws.Range("D41").Value = '2.62'
$ws.Range("E41").Value = '  +2.96%  '
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D44").Value = '119.86'
$ws.Range("E44").Value = '  -4.50%  '
$ws.Range("D45").Value = '21.82'
$ws.Range("E45").Value = '  -7.46%  '
$ws.Range("D46").Value = '2.083.80'
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("E47").Value = '  -4.04%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '5.75'
$ws.Range("E49").Value = '  -3.51%  '
$ws.Range("B50").Value = 'SEI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D50").Value = '0.956'
$ws.Range("E50").Value = '  -3.43%  '
$ws.Range("D51").Value = '1.15'
$ws.Range("E51").Value = '  +32.63%  '
